$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.965.02'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.556.65'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.32'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.12'
$ws.Range('E8').Value = '  +4.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '1.778.82'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '1.556.81'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').Value = '26.970.28'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.80'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '218.38'
$ws.Range('E18').Value = '  +2.43%  '
$ws.Range('D19').Value = '0.0₃0697'
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.31'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.16'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.99'
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +2.55%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('D33').Value = '1.423.64'
$ws.Range('E33').Value = '  +5.05%  '
$ws.Range('E34').Value = '  +5.20%  '
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.981'
$ws.Range('E36').Value = '  +1.66%  '
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.521'
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('E43').Value = '  +4.97%  '
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.49'
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '1.692.28'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.04'
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('D50').Value = '0.0₇0999'
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('E51').Value = '  +1.09%  '
